$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 17858.334
$ws.Range("I7").Value = 2383.3333
$ws.Range("K7").Value = 2383.3333
$ws.Range("M7").Value = -2271.3333
$ws.Range("H10").Value = 2002
$ws.Range("I10").Value = 2002
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 2002
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -1709
$ws.Range("H14").Value = 17858.334
$ws.Range("I14").Value = 2383.3333
$ws.Range("K14").Value = 2383.3333
$ws.Range("M14").Value = -2192.3333
$ws.Range("H28").Value = 64158.875
$ws.Range("I28").Value = 78441.234
$ws.Range("J28").Value = 2268.6667
$ws.Range("K28").Value = 78441.234
$ws.Range("L28").Value = 2268.6667
$ws.Range("M28").Value = -77956.234
$ws.Range("N28").Value = -3238.6667
$ws.Range("H41").Value = 1410.6666
$ws.Range("I41").Value = 1137.6154
$ws.Range("J41").Value = 2120.6
$ws.Range("K41").Value = 1137.6154
$ws.Range("L41").Value = 2120.6
$ws.Range("M41").Value = -697.6153999999999
$ws.Range("N41").Value = -3000.6
$ws.Range("H62").Value = 15197.363
$ws.Range("I62").Value = 18741.691
$ws.Range("J62").Value = 10077.777
$ws.Range("K62").Value = 18741.691
$ws.Range("L62").Value = 10077.777
$ws.Range("M62").Value = -18117.691
$ws.Range("N62").Value = -11325.777
$ws.Range("H65").Value = 15197.363
$ws.Range("I65").Value = 18741.691
$ws.Range("J65").Value = 10077.777
$ws.Range("K65").Value = 93708.45499999999
$ws.Range("L65").Value = 50388.885
$ws.Range("M65").Value = -90588.45499999999
$ws.Range("N65").Value = -56628.885
$ws.Range("H107").Value = 43478616
$ws.Range("I107").Value = 55555856
$ws.Range("J107").Value = 555
$ws.Range("K107").Value = 55555856
$ws.Range("L107").Value = 555
$ws.Range("M107").Value = -55553936
$ws.Range("N107").Value = -4395
$ws.Range("H116").Value = 82492.8
$ws.Range("I116").Value = 133321.5
$ws.Range("J116").Value = 6249.75
$ws.Range("K116").Value = 133321.5
$ws.Range("L116").Value = 6249.75
$ws.Range("M116").Value = -129879.5
$ws.Range("N116").Value = -13133.75
$ws.Range("H133").Value = 105248.625
$ws.Range("J133").Value = 105248.625
$ws.Range("L133").Value = 105248.625
$ws.Range("N133").Value = -115368.625
$ws.Range("H136").Value = 99994
$ws.Range("J136").Value = 99994
$ws.Range("L136").Value = 99994
$ws.Range("N136").Value = -110194
$ws.Range("H137").Value = 2502.5
$ws.Range("I137").Value = 1071.2778
$ws.Range("K137").Value = 3213.8334
$ws.Range("M137").Value = -663.8334000000004
$ws.Range("H138").Value = 3454.8965
$ws.Range("I138").Value = 1466.5454
$ws.Range("J138").Value = 4670
$ws.Range("K138").Value = 4399.6362
$ws.Range("L138").Value = 14010
$ws.Range("M138").Value = 740.3638000000001
$ws.Range("N138").Value = -24290

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 934.43396
$ws.Range("I32").Value = 934.43396
$ws.Range("K32").Value = 934.43396
$ws.Range("M32").Value = -647.43396
$ws.Range("H35").Value = 4666.6665
$ws.Range("I35").Value = 4666.6665
$ws.Range("K35").Value = 4666.6665
$ws.Range("M35").Value = -4260.6665
$ws.Range("H61").Value = 2759.8635
$ws.Range("I61").Value = 2643.7896
$ws.Range("K61").Value = 2643.7896
$ws.Range("M61").Value = -2431.7896
$ws.Range("H131").Value = 54975
$ws.Range("J131").Value = 54975
$ws.Range("L131").Value = 54975
$ws.Range("N131").Value = -65055
$ws.Range("H132").Value = 2918.3333
$ws.Range("I132").Value = 2839.275
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 8517.825000000001
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -5987.825000000001
$ws.Range("N132").Value = -18558.5
$ws.Range("H136").Value = 2759.8635
$ws.Range("I136").Value = 2643.7896
$ws.Range("K136").Value = 7931.3688
$ws.Range("M136").Value = -5381.3688

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2572.0571
$ws.Range("I105").Value = 2641.0356
$ws.Range("K105").Value = 2641.0356
$ws.Range("M105").Value = -894.0356000000002
$ws.Range("H107").Value = 1553.5714
$ws.Range("I107").Value = 1946.3334
$ws.Range("J107").Value = 1071.5454
$ws.Range("K107").Value = 1946.3334
$ws.Range("L107").Value = 1071.5454
$ws.Range("M107").Value = -26.33339999999998
$ws.Range("N107").Value = -4911.5454
$ws.Range("H134").Value = 11388.48
$ws.Range("I134").Value = 12192.714
$ws.Range("J134").Value = 7166.25
$ws.Range("K134").Value = 36578.142
$ws.Range("L134").Value = 21498.75
$ws.Range("M134").Value = -34043.142
$ws.Range("N134").Value = -26568.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 11800
$ws.Range("I33").Value = 10500
$ws.Range("J33").Value = 17000
$ws.Range("K33").Value = 10500
$ws.Range("L33").Value = 17000
$ws.Range("M33").Value = -10121
$ws.Range("N33").Value = -17758
$ws.Range("H62").Value = 3386.875
$ws.Range("I62").Value = 3199.25
$ws.Range("J62").Value = 3574.5
$ws.Range("K62").Value = 3199.25
$ws.Range("L62").Value = 3574.5
$ws.Range("M62").Value = -2575.25
$ws.Range("N62").Value = -4822.5
$ws.Range("H65").Value = 3386.875
$ws.Range("I65").Value = 3199.25
$ws.Range("J65").Value = 3574.5
$ws.Range("K65").Value = 15996.25
$ws.Range("L65").Value = 17872.5
$ws.Range("M65").Value = -12876.25
$ws.Range("N65").Value = -24112.5
$ws.Range("H105").Value = 57985.93
$ws.Range("J105").Value = 1291
$ws.Range("L105").Value = 1291
$ws.Range("N105").Value = -4785
$ws.Range("H107").Value = 705.7368
$ws.Range("I107").Value = 626.6923
$ws.Range("J107").Value = 877
$ws.Range("K107").Value = 626.6923
$ws.Range("L107").Value = 877
$ws.Range("M107").Value = 1293.3077
$ws.Range("N107").Value = -4717
$ws.Range("H114").Value = 88012.75
$ws.Range("J114").Value = 88012.75
$ws.Range("L114").Value = 88012.75
$ws.Range("N114").Value = -96690.75
$ws.Range("H132").Value = 4338.1113
$ws.Range("I132").Value = 4503.875
$ws.Range("J132").Value = 3012
$ws.Range("K132").Value = 13511.625
$ws.Range("L132").Value = 9036
$ws.Range("M132").Value = -10981.625
$ws.Range("N132").Value = -14096

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 471.2857
$ws.Range("I19").Value = 275
$ws.Range("J19").Value = 549.8
$ws.Range("K19").Value = 825
$ws.Range("L19").Value = 1649.4
$ws.Range("M19").Value = -651
$ws.Range("N19").Value = -1997.4
$ws.Range("H23").Value = 811.44446
$ws.Range("I23").Value = 588
$ws.Range("J23").Value = 839.375
$ws.Range("K23").Value = 1764
$ws.Range("L23").Value = 2518.125
$ws.Range("M23").Value = -1529
$ws.Range("N23").Value = -2988.125
$ws.Range("H34").Value = 4346.6
$ws.Range("I34").Value = 244.33333
$ws.Range("J34").Value = 10500
$ws.Range("K34").Value = 732.99999
$ws.Range("L34").Value = 31500
$ws.Range("M34").Value = -648.99999
$ws.Range("N34").Value = -31668
$ws.Range("H86").Value = 696
$ws.Range("I86").Value = 927.25
$ws.Range("J86").Value = 511
$ws.Range("K86").Value = 2781.75
$ws.Range("L86").Value = 1533
$ws.Range("M86").Value = -1595.75
$ws.Range("N86").Value = -3905
$ws.Range("H89").Value = 696
$ws.Range("I89").Value = 927.25
$ws.Range("J89").Value = 511
$ws.Range("K89").Value = 8345.25
$ws.Range("L89").Value = 4599
$ws.Range("M89").Value = -2417.25
$ws.Range("N89").Value = -16455
$ws.Range("H93").Value = 5862
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 5862
$ws.Range("K93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("M93").Value = 17586
$ws.Range("N93").Value = -21330
$ws.Range("H137").Value = 4840.6665
$ws.Range("I137").Value = 800
$ws.Range("J137").Value = 5208
$ws.Range("K137").Value = 2400
$ws.Range("L137").Value = 15624
$ws.Range("M137").Value = 2700
$ws.Range("N137").Value = -25824

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3516.5
$ws.Range("I113").Value = 2421.2727
$ws.Range("J113").Value = 4855.1113
$ws.Range("K113").Value = 2421.2727
$ws.Range("L113").Value = 4855.1113
$ws.Range("M113").Value = -251.2727
$ws.Range("N113").Value = -9195.1113
$ws.Range("H132").Value = 6682.375
$ws.Range("I132").Value = 6194.533
$ws.Range("K132").Value = 18583.599
$ws.Range("M132").Value = -16053.599

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2188.9
$ws.Range("J61").Value = 2839.8
$ws.Range("L61").Value = 2839.8
$ws.Range("N61").Value = -3243.8
$ws.Range("H68").Value = 2248.4736
$ws.Range("I68").Value = 1916.1428
$ws.Range("J68").Value = 3179
$ws.Range("K68").Value = 1916.1428
$ws.Range("L68").Value = 3179
$ws.Range("M68").Value = -1167.1428
$ws.Range("N68").Value = -4677
$ws.Range("H71").Value = 2248.4736
$ws.Range("I71").Value = 1916.1428
$ws.Range("J71").Value = 3179
$ws.Range("K71").Value = 9580.714
$ws.Range("L71").Value = 15895
$ws.Range("M71").Value = -5836.714
$ws.Range("N71").Value = -23383
$ws.Range("H113").Value = 2188.9
$ws.Range("J113").Value = 2839.8
$ws.Range("L113").Value = 2839.8
$ws.Range("N113").Value = -7179.8
$ws.Range("H122").Value = 51194.137
$ws.Range("I122").Value = 5484.6
$ws.Range("J122").Value = 149143.14
$ws.Range("K122").Value = 16453.8
$ws.Range("L122").Value = 447429.42
$ws.Range("M122").Value = -14003.8
$ws.Range("N122").Value = -452329.42
$ws.Range("H132").Value = 17749.285
$ws.Range("I132").Value = 18451.8
$ws.Range("K132").Value = 55355.39999999999
$ws.Range("M132").Value = -52825.39999999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1392.1538
$ws.Range("I122").Value = 1341.7084
$ws.Range("J122").Value = 1997.5
$ws.Range("K122").Value = 4025.1252
$ws.Range("L122").Value = 5992.5
$ws.Range("M122").Value = -10892.5
